# Auto-generated edit script: apply numeric updates to Behemoth_Profits (Leve Profit tables)
# across all 8 job-abbreviation sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1277.3846
$ws.Range("I41").Value = 862.8
$ws.Range("J41").Value = 1536.5
$ws.Range("K41").Value = 862.8
$ws.Range("L41").Value = 1536.5
$ws.Range("M41").Value = -422.8
$ws.Range("N41").Value = -2416.5
$ws.Range("H69").Value = 11926.375
$ws.Range("J69").Value = 12915.857
$ws.Range("L69").Value = 38747.571
$ws.Range("N69").Value = -40495.571
$ws.Range("H70").Value = 2507.5
$ws.Range("J70").Value = 2507.5
$ws.Range("L70").Value = 7522.5
$ws.Range("N70").Value = -8062.5
$ws.Range("H72").Value = 11926.375
$ws.Range("J72").Value = 12915.857
$ws.Range("L72").Value = 116242.713
$ws.Range("N72").Value = -124978.713
$ws.Range("H73").Value = 2507.5
$ws.Range("J73").Value = 2507.5
$ws.Range("L73").Value = 7522.5
$ws.Range("N73").Value = -9394.5
$ws.Range("H81").Value = 46584.75
$ws.Range("J81").Value = 46584.75
$ws.Range("L81").Value = 46584.75
$ws.Range("N81").Value = -48580.75
$ws.Range("H82").Value = 35294.5
$ws.Range("J82").Value = 70047
$ws.Range("L82").Value = 210141
$ws.Range("N82").Value = -210953
$ws.Range("H84").Value = 46584.75
$ws.Range("J84").Value = 46584.75
$ws.Range("L84").Value = 139754.25
$ws.Range("N84").Value = -149738.25
$ws.Range("H85").Value = 35294.5
$ws.Range("J85").Value = 70047
$ws.Range("L85").Value = 210141
$ws.Range("N85").Value = -212949
$ws.Range("H87").Value = 119996.664
$ws.Range("J87").Value = 119996.664
$ws.Range("L87").Value = 119996.664
$ws.Range("N87").Value = -122492.664
$ws.Range("H90").Value = 119996.664
$ws.Range("J90").Value = 119996.664
$ws.Range("L90").Value = 359989.992
$ws.Range("N90").Value = -372469.992
$ws.Range("H92").Value = 1030.138
$ws.Range("I92").Value = 912.5
$ws.Range("K92").Value = 912.5
$ws.Range("M92").Value = 335.5
$ws.Range("H93").Value = 53977.535
$ws.Range("J93").Value = 55690.215
$ws.Range("L93").Value = 55690.215
$ws.Range("N93").Value = -60682.215
$ws.Range("H95").Value = 68244.5
$ws.Range("J95").Value = 68244.5
$ws.Range("L95").Value = 68244.5
$ws.Range("N95").Value = -73736.5
$ws.Range("H137").Value = 4019.7036
$ws.Range("I137").Value = 2554.8235
$ws.Range("K137").Value = 7664.470499999999
$ws.Range("M137").Value = -5114.470499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35717948
$ws.Range("I32").Value = 45457480
$ws.Range("J32").Value = 6333
$ws.Range("K32").Value = 45457480
$ws.Range("L32").Value = 6333
$ws.Range("M32").Value = -45457193
$ws.Range("N32").Value = -6907
$ws.Range("H61").Value = 35791384
$ws.Range("I61").Value = 71430750
$ws.Range("K61").Value = 71430750
$ws.Range("M61").Value = -71430538
$ws.Range("H88").Value = 1521.7646
$ws.Range("I88").Value = 1467.375
$ws.Range("J88").Value = 1570.1111
$ws.Range("K88").Value = 1467.375
$ws.Range("L88").Value = 1570.1111
$ws.Range("M88").Value = -1061.375
$ws.Range("N88").Value = -2382.1111
$ws.Range("H91").Value = 1521.7646
$ws.Range("I91").Value = 1467.375
$ws.Range("J91").Value = 1570.1111
$ws.Range("K91").Value = 1467.375
$ws.Range("L91").Value = 1570.1111
$ws.Range("M91").Value = -63.375
$ws.Range("N91").Value = -4378.1111
$ws.Range("H94").Value = 43721.8
$ws.Range("J94").Value = 43721.8
$ws.Range("L94").Value = 43721.8
$ws.Range("N94").Value = -45523.8
$ws.Range("H106").Value = 21943.8
$ws.Range("J106").Value = 21943.8
$ws.Range("L106").Value = 21943.8
$ws.Range("N106").Value = -24467.8
$ws.Range("H122").Value = 1811.1
$ws.Range("I122").Value = 1456.8889
$ws.Range("K122").Value = 4370.6667
$ws.Range("M122").Value = -1920.6667
$ws.Range("H134").Value = 290000
$ws.Range("J134").Value = 290000
$ws.Range("L134").Value = 290000
$ws.Range("N134").Value = -300140
$ws.Range("H136").Value = 35791384
$ws.Range("I136").Value = 71430750
$ws.Range("K136").Value = 214292250
$ws.Range("M136").Value = -214289700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5346.722
$ws.Range("I20").Value = 4662.1333
$ws.Range("J20").Value = 8769.666999999999
$ws.Range("K20").Value = 4662.1333
$ws.Range("L20").Value = 8769.666999999999
$ws.Range("M20").Value = -4415.1333
$ws.Range("N20").Value = -9263.666999999999
$ws.Range("H40").Value = 81132.60000000001
$ws.Range("J40").Value = 77744.664
$ws.Range("L40").Value = 77744.664
$ws.Range("N40").Value = -78274.664
$ws.Range("H96").Value = 34747.4
$ws.Range("I96").Value = 15999.8
$ws.Range("J96").Value = 53495
$ws.Range("K96").Value = 15999.8
$ws.Range("L96").Value = 53495
$ws.Range("M96").Value = -13253.8
$ws.Range("N96").Value = -58987

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1369816.8
$ws.Range("I31").Value = 1219.7858
$ws.Range("J31").Value = 7756602.5
$ws.Range("K31").Value = 1219.7858
$ws.Range("L31").Value = 7756602.5
$ws.Range("M31").Value = -924.7858000000001
$ws.Range("N31").Value = -7757192.5
$ws.Range("H34").Value = 1369816.8
$ws.Range("I34").Value = 1219.7858
$ws.Range("J34").Value = 7756602.5
$ws.Range("K34").Value = 1219.7858
$ws.Range("L34").Value = 7756602.5
$ws.Range("M34").Value = -1017.7858
$ws.Range("N34").Value = -7757006.5
$ws.Range("H62").Value = 2975
$ws.Range("J62").Value = 3166.6667
$ws.Range("L62").Value = 3166.6667
$ws.Range("N62").Value = -4414.6667
$ws.Range("H65").Value = 2975
$ws.Range("J65").Value = 3166.6667
$ws.Range("L65").Value = 15833.3335
$ws.Range("N65").Value = -22073.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H82").Value = 9996.666999999999
$ws.Range("J82").Value = 9996.666999999999
$ws.Range("L82").Value = 29990.001
$ws.Range("N82").Value = -30802.001
$ws.Range("H85").Value = 9996.666999999999
$ws.Range("J85").Value = 9996.666999999999
$ws.Range("L85").Value = 29990.001
$ws.Range("N85").Value = -32798.001
$ws.Range("H113").Value = 999.46155
$ws.Range("I113").Value = 408.33334
$ws.Range("J113").Value = 1176.8
$ws.Range("K113").Value = 1225.00002
$ws.Range("L113").Value = 3530.4
$ws.Range("M113").Value = 944.9999800000001
$ws.Range("N113").Value = -7870.4
$ws.Range("H127").Value = 2999.0557
$ws.Range("J127").Value = 2999.0557
$ws.Range("L127").Value = 8997.167099999999
$ws.Range("N127").Value = -18917.1671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 116000
$ws.Range("J64").Value = 116000
$ws.Range("L64").Value = 116000
$ws.Range("N64").Value = -116496
$ws.Range("H67").Value = 116000
$ws.Range("J67").Value = 116000
$ws.Range("L67").Value = 116000
$ws.Range("N67").Value = -117716
$ws.Range("H105").Value = 108412.5
$ws.Range("J105").Value = 108412.5
$ws.Range("L105").Value = 108412.5
$ws.Range("N105").Value = -115400.5
$ws.Range("H126").Value = 4699
$ws.Range("I126").Value = 4531.778
$ws.Range("K126").Value = 13595.334
$ws.Range("M126").Value = -11125.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 67031.69
$ws.Range("I7").Value = 3873.9167
$ws.Range("J7").Value = 256505
$ws.Range("K7").Value = 3873.9167
$ws.Range("L7").Value = 256505
$ws.Range("M7").Value = -3761.9167
$ws.Range("N7").Value = -256729
$ws.Range("H16").Value = 1808.6
$ws.Range("I16").Value = 1808.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1808.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1638.6
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 3798.611
$ws.Range("I40").Value = 2952
$ws.Range("K40").Value = 2952
$ws.Range("M40").Value = -2816
$ws.Range("H46").Value = 3193.8635
$ws.Range("I46").Value = 3124.0908
$ws.Range("K46").Value = 3124.0908
$ws.Range("M46").Value = -2936.0908
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 62500210
$ws.Range("I55").Value = 76923280
$ws.Range("J55").Value = 240.66667
$ws.Range("K55").Value = 76923280
$ws.Range("L55").Value = 240.66667
$ws.Range("M55").Value = -76923107
$ws.Range("N55").Value = -586.6666700000001
$ws.Range("H103").Value = 34915.168
$ws.Range("J103").Value = 34915.168
$ws.Range("L103").Value = 34915.168
$ws.Range("N103").Value = -37259.168
$ws.Range("H126").Value = 67031.69
$ws.Range("I126").Value = 3873.9167
$ws.Range("J126").Value = 256505
$ws.Range("K126").Value = 11621.7501
$ws.Range("L126").Value = 769515
$ws.Range("M126").Value = -9151.750100000001
$ws.Range("N126").Value = -774455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
